# Force Sensor Validation - update input voltage and R2 resistance values
# M1: Input Voltage 5 -> 3.3
# K3:K41: R2 Resistance 3300 -> 10000 (L/M columns recalc automatically)
# Move the active selection to P5

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("M1").Value = 3.3
$ws.Range("K3:K41").Value = 10000

$ws.Range("P5").Select()
